$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (matching the
# original inlineStr cell type) instead of letting Excel auto-convert it to a
# number. We flag the cell as Text, assign the value, then restore the "Normal"
# cell style so no stray formatting is left behind on the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "42.630.16"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "2.364.83"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  -0.13%  "

Set-TextValue "D5" "329.06"
$ws.Range("E5").Value = "  +5.70%  "

Set-TextValue "D6" "100.75"
$ws.Range("E6").Value = "  -7.42%  "

Set-TextValue "D7" "0.637"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  +0.18%  "

Set-TextValue "D9" "0.618"
$ws.Range("E9").Value = "  -0.96%  "

Set-TextValue "D10" "40.26"
$ws.Range("E10").Value = "  -7.65%  "

Set-TextValue "D11" "0.0921"
$ws.Range("E11").Value = "  -1.75%  "

Set-TextValue "D12" "8.40"
$ws.Range("E12").Value = "  -6.15%  "

Set-TextValue "D13" "1.02"
$ws.Range("E13").Value = "  -6.02%  "

$ws.Range("E14").Value = "  +0.41%  "

Set-TextValue "D15" "16.31"
$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").Value = "2.724.92"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "2.360.48"
$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").Value = "42.592.22"
$ws.Range("E18").Value = "  -0.99%  "

Set-TextValue "D19" "7.72"
$ws.Range("E19").Value = "  +6.77%  "

$ws.Range("E20").Value = "  -1.85%  "

Set-TextValue "D21" "75.28"
$ws.Range("E21").Value = "  -0.56%  "

Set-TextValue "D22" "3.69"
$ws.Range("E22").Value = "  +6.56%  "

$ws.Range("E23").Value = "  +10.99%  "

$ws.Range("E24").Value = "  -8.32%  "

Set-TextValue "D25" "9.72"
$ws.Range("E25").Value = "  +8.41%  "

Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.13%  "

Set-TextValue "D27" "11.48"
$ws.Range("E27").Value = "  -3.69%  "

Set-TextValue "D28" "23.91"
$ws.Range("E28").Value = "  +5.99%  "

$ws.Range("E29").Value = "  -1.49%  "

Set-TextValue "D30" "173.90"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("E31").Value = "  -1.89%  "

$ws.Range("E32").Value = "  -1.05%  "

Set-TextValue "D33" "35.47"
$ws.Range("E33").Value = "  -8.55%  "

Set-TextValue "D34" "5.96"
$ws.Range("E34").Value = "  +2.39%  "

$ws.Range("E35").Value = "  +1.25%  "

Set-TextValue "D36" "4.60"
$ws.Range("E36").Value = "  -7.65%  "

Set-TextValue "D37" "0.0359"
$ws.Range("E37").Value = "  -5.19%  "

$ws.Range("E38").Value = "  -5.43%  "

$ws.Range("E39").Value = "  +4.14%  "

Set-TextValue "D40" "0.105"
$ws.Range("E40").Value = "  +0.70%  "

Set-TextValue "D41" "1.52"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("E42").Value = "  -2.57%  "

Set-TextValue "D43" "68.76"
$ws.Range("E43").Value = "  -4.66%  "

$ws.Range("E44").Value = "  +0.03%  "

Set-TextValue "D45" "115.08"
$ws.Range("E45").Value = "  +4.18%  "

$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D46" "88.31"
$ws.Range("E46").Value = "  +38.19%  "

$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D47" "11.96"
$ws.Range("E47").Value = "  -4.54%  "

Set-TextValue "D48" "5.46"
$ws.Range("E48").Value = "  -4.33%  "

Set-TextValue "D49" "8.99"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("D50").Value = "1.593.02"
$ws.Range("E50").Value = "  +6.77%  "

$ws.Range("E51").Value = "  -3.28%  "
